# issue #5: stock data from json to db
# Sheet "股票" (stock, sheet index 5) gains three new columns:
#   - "category"    inserted right after "property_category" (between old H and I)
#   - "source_file" inserted right after "legislator_id" (which itself shifts right)
#   - "index"       inserted right after "source_file", duplicating the row id in col A
# Also fixes a stray leading glyph in one quantity value ("<80>100" -> "100").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Insert a single new column I (old date/legislator_name/legislator_id shift to J/K/L).
$ws.Range("I1:I6").EntireColumn.Insert()

# Insert two new columns at M:N (after the now-shifted legislator_id in column L).
$ws.Range("M1:N6").EntireColumn.Insert()

# --- Header row ---
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Data rows: category / source_file / index ---
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmpc7221"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value
}

# --- Fix stray leading glyph on row 6's quantity text ---
$ws.Cells.Item(6, 4).Value = "100"
